$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = "asdvdv"
$ws.Range("C3").Value = "sdv"
$ws.Range("B4").Value = "33r"

$ws.Range("B4").Select()
